$d = $word.ActiveDocument

$d.Content.Find.Execute("7+6=13", $true, $false, $false, $false, $false, $true, 1, $false, "71+28=99", 2) | Out-Null
$d.Content.Find.Execute("36+2=38", $true, $false, $false, $false, $false, $true, 1, $false, "56-0=56", 2) | Out-Null
$d.Content.Find.Execute("63-53=10", $true, $false, $false, $false, $false, $true, 1, $false, "27+54=81", 2) | Out-Null
$d.Content.Find.Execute("47+6=53", $true, $false, $false, $false, $false, $true, 1, $false, "94-61=33", 2) | Out-Null
$d.Content.Find.Execute("78+4=82", $true, $false, $false, $false, $false, $true, 1, $false, "67-9=58", 2) | Out-Null
$d.Content.Find.Execute("16+23=39", $true, $false, $false, $false, $false, $true, 1, $false, "40-9=31", 2) | Out-Null
$d.Content.Find.Execute("44+13=57", $true, $false, $false, $false, $false, $true, 1, $false, "75+10=85", 2) | Out-Null
$d.Content.Find.Execute("49-17=32", $true, $false, $false, $false, $false, $true, 1, $false, "0+27=27", 2) | Out-Null
$d.Content.Find.Execute("38+40=78", $true, $false, $false, $false, $false, $true, 1, $false, "5+54=59", 2) | Out-Null
$d.Content.Find.Execute("42-21=21", $true, $false, $false, $false, $false, $true, 1, $false, "1+60=61", 2) | Out-Null
$d.Content.Find.Execute("7+4=11", $true, $false, $false, $false, $false, $true, 1, $false, "47-10=37", 2) | Out-Null
$d.Content.Find.Execute("85-23=62", $true, $false, $false, $false, $false, $true, 1, $false, "28+35=63", 2) | Out-Null
$d.Content.Find.Execute("21+1=22", $true, $false, $false, $false, $false, $true, 1, $false, "93-29=64", 2) | Out-Null
$d.Content.Find.Execute("83+8=91", $true, $false, $false, $false, $false, $true, 1, $false, "50+49=99", 2) | Out-Null
$d.Content.Find.Execute("17+81=98", $true, $false, $false, $false, $false, $true, 1, $false, "9+77=86", 2) | Out-Null
$d.Content.Find.Execute("20+30=50", $true, $false, $false, $false, $false, $true, 1, $false, "97-89=8", 2) | Out-Null
$d.Content.Find.Execute("16+36=52", $true, $false, $false, $false, $false, $true, 1, $false, "21+0=21", 2) | Out-Null
$d.Content.Find.Execute("71-58=13", $true, $false, $false, $false, $false, $true, 1, $false, "53+31=84", 2) | Out-Null
$d.Content.Find.Execute("88-50=38", $true, $false, $false, $false, $false, $true, 1, $false, "3+65=68", 2) | Out-Null
$d.Content.Find.Execute("43+56=99", $true, $false, $false, $false, $false, $true, 1, $false, "0+95=95", 2) | Out-Null
$d.Content.Find.Execute("13+69=82", $true, $false, $false, $false, $false, $true, 1, $false, "36+30=66", 2) | Out-Null
$d.Content.Find.Execute("30-23=7", $true, $false, $false, $false, $false, $true, 1, $false, "28+8=36", 2) | Out-Null
$d.Content.Find.Execute("3+54=57", $true, $false, $false, $false, $false, $true, 1, $false, "35-4=31", 2) | Out-Null
$d.Content.Find.Execute("60+7=67", $true, $false, $false, $false, $false, $true, 1, $false, "75+18=93", 2) | Out-Null
$d.Content.Find.Execute("55+12=67", $true, $false, $false, $false, $false, $true, 1, $false, "61-13=48", 2) | Out-Null
$d.Content.Find.Execute("63-19=44", $true, $false, $false, $false, $false, $true, 1, $false, "74-64=10", 2) | Out-Null
$d.Content.Find.Execute("31+40=71", $true, $false, $false, $false, $false, $true, 1, $false, "30-26=4", 2) | Out-Null
$d.Content.Find.Execute("16+53=69", $true, $false, $false, $false, $false, $true, 1, $false, "11+9=20", 2) | Out-Null
$d.Content.Find.Execute("91+7=98", $true, $false, $false, $false, $false, $true, 1, $false, "20+54=74", 2) | Out-Null
$d.Content.Find.Execute("41-11=30", $true, $false, $false, $false, $false, $true, 1, $false, "63-24=39", 2) | Out-Null
$d.Content.Find.Execute("77-11=66", $true, $false, $false, $false, $false, $true, 1, $false, "51-33=18", 2) | Out-Null
$d.Content.Find.Execute("74-5=69", $true, $false, $false, $false, $false, $true, 1, $false, "46+47=93", 2) | Out-Null
$d.Content.Find.Execute("90-15=75", $true, $false, $false, $false, $false, $true, 1, $false, "38+18=56", 2) | Out-Null
$d.Content.Find.Execute("50+47=97", $true, $false, $false, $false, $false, $true, 1, $false, "66+2=68", 2) | Out-Null
$d.Content.Find.Execute("46-34=12", $true, $false, $false, $false, $false, $true, 1, $false, "50+29=79", 2) | Out-Null
$d.Content.Find.Execute("1+10=11", $true, $false, $false, $false, $false, $true, 1, $false, "24+24=48", 2) | Out-Null
$d.Content.Find.Execute("93-93=0", $true, $false, $false, $false, $false, $true, 1, $false, "79-6=73", 2) | Out-Null
$d.Content.Find.Execute("18+67=85", $true, $false, $false, $false, $false, $true, 1, $false, "93-93=0", 2) | Out-Null
$d.Content.Find.Execute("19+71=90", $true, $false, $false, $false, $false, $true, 1, $false, "24+71=95", 2) | Out-Null
$d.Content.Find.Execute("11+17=28", $true, $false, $false, $false, $false, $true, 1, $false, "25+8=33", 2) | Out-Null
$d.Content.Find.Execute("57-39=18", $true, $false, $false, $false, $false, $true, 1, $false, "81-2=79", 2) | Out-Null
$d.Content.Find.Execute("29+36=65", $true, $false, $false, $false, $false, $true, 1, $false, "28-23=5", 2) | Out-Null
$d.Content.Find.Execute("89-24=65", $true, $false, $false, $false, $false, $true, 1, $false, "69-61=8", 2) | Out-Null
$d.Content.Find.Execute("7-6=1", $true, $false, $false, $false, $false, $true, 1, $false, "2+77=79", 2) | Out-Null
$d.Content.Find.Execute("69-65=4", $true, $false, $false, $false, $false, $true, 1, $false, "41+27=68", 2) | Out-Null
$d.Content.Find.Execute("63-35=28", $true, $false, $false, $false, $false, $true, 1, $false, "10-4=6", 2) | Out-Null
$d.Content.Find.Execute("28+38=66", $true, $false, $false, $false, $false, $true, 1, $false, "3+82=85", 2) | Out-Null
$d.Content.Find.Execute("31-17=14", $true, $false, $false, $false, $false, $true, 1, $false, "52+9=61", 2) | Out-Null
$d.Content.Find.Execute("41+41=82", $true, $false, $false, $false, $false, $true, 1, $false, "37+56=93", 2) | Out-Null
$d.Content.Find.Execute("98-74=24", $true, $false, $false, $false, $false, $true, 1, $false, "85-73=12", 2) | Out-Null
$d.Content.Find.Execute("84-29=55", $true, $false, $false, $false, $false, $true, 1, $false, "97-60=37", 2) | Out-Null
$d.Content.Find.Execute("21+38=59", $true, $false, $false, $false, $false, $true, 1, $false, "26+24=50", 2) | Out-Null
$d.Content.Find.Execute("76+8=84", $true, $false, $false, $false, $false, $true, 1, $false, "97-21=76", 2) | Out-Null
$d.Content.Find.Execute("21+55=76", $true, $false, $false, $false, $false, $true, 1, $false, "2+28=30", 2) | Out-Null
$d.Content.Find.Execute("65+34=99", $true, $false, $false, $false, $false, $true, 1, $false, "1+15=16", 2) | Out-Null
$d.Content.Find.Execute("57-49=8", $true, $false, $false, $false, $false, $true, 1, $false, "27+22=49", 2) | Out-Null
$d.Content.Find.Execute("72+25=97", $true, $false, $false, $false, $false, $true, 1, $false, "89-65=24", 2) | Out-Null
$d.Content.Find.Execute("76-2=74", $true, $false, $false, $false, $false, $true, 1, $false, "17+45=62", 2) | Out-Null
$d.Content.Find.Execute("23-20=3", $true, $false, $false, $false, $false, $true, 1, $false, "33-18=15", 2) | Out-Null
$d.Content.Find.Execute("61+16=77", $true, $false, $false, $false, $false, $true, 1, $false, "69-64=5", 2) | Out-Null
$d.Content.Find.Execute("20+79=99", $true, $false, $false, $false, $false, $true, 1, $false, "25+32=57", 2) | Out-Null
$d.Content.Find.Execute("42+35=77", $true, $false, $false, $false, $false, $true, 1, $false, "89-3=86", 2) | Out-Null
$d.Content.Find.Execute("13+20=33", $true, $false, $false, $false, $false, $true, 1, $false, "93-17=76", 2) | Out-Null
$d.Content.Find.Execute("89-30=59", $true, $false, $false, $false, $false, $true, 1, $false, "55+24=79", 2) | Out-Null
$d.Content.Find.Execute("69-31=38", $true, $false, $false, $false, $false, $true, 1, $false, "91-84=7", 2) | Out-Null
$d.Content.Find.Execute("45-41=4", $true, $false, $false, $false, $false, $true, 1, $false, "89-65=24", 2) | Out-Null
$d.Content.Find.Execute("44+42=86", $true, $false, $false, $false, $false, $true, 1, $false, "36-7=29", 2) | Out-Null
$d.Content.Find.Execute("98-43=55", $true, $false, $false, $false, $false, $true, 1, $false, "3+47=50", 2) | Out-Null
$d.Content.Find.Execute("6+54=60", $true, $false, $false, $false, $false, $true, 1, $false, "3+3=6", 2) | Out-Null
$d.Content.Find.Execute("38+22=60", $true, $false, $false, $false, $false, $true, 1, $false, "21+7=28", 2) | Out-Null
$d.Content.Find.Execute("89-47=42", $true, $false, $false, $false, $false, $true, 1, $false, "68-63=5", 2) | Out-Null
$d.Content.Find.Execute("29+4=33", $true, $false, $false, $false, $false, $true, 1, $false, "49-22=27", 2) | Out-Null
$d.Content.Find.Execute("76+1=77", $true, $false, $false, $false, $false, $true, 1, $false, "65+28=93", 2) | Out-Null
$d.Content.Find.Execute("53-48=5", $true, $false, $false, $false, $false, $true, 1, $false, "99-79=20", 2) | Out-Null
$d.Content.Find.Execute("15+25=40", $true, $false, $false, $false, $false, $true, 1, $false, "10-7=3", 2) | Out-Null
$d.Content.Find.Execute("0+37=37", $true, $false, $false, $false, $false, $true, 1, $false, "87-64=23", 2) | Out-Null
$d.Content.Find.Execute("88-27=61", $true, $false, $false, $false, $false, $true, 1, $false, "70-30=40", 2) | Out-Null
$d.Content.Find.Execute("39+39=78", $true, $false, $false, $false, $false, $true, 1, $false, "34+27=61", 2) | Out-Null
$d.Content.Find.Execute("48-33=15", $true, $false, $false, $false, $false, $true, 1, $false, "19+15=34", 2) | Out-Null
$d.Content.Find.Execute("70+19=89", $true, $false, $false, $false, $false, $true, 1, $false, "50+40=90", 2) | Out-Null
$d.Content.Find.Execute("58+18=76", $true, $false, $false, $false, $false, $true, 1, $false, "36+22=58", 2) | Out-Null
$d.Content.Find.Execute("31-5=26", $true, $false, $false, $false, $false, $true, 1, $false, "94-66=28", 2) | Out-Null
$d.Content.Find.Execute("33-25=8", $true, $false, $false, $false, $false, $true, 1, $false, "90-86=4", 2) | Out-Null
$d.Content.Find.Execute("14+82=96", $true, $false, $false, $false, $false, $true, 1, $false, "99-92=7", 2) | Out-Null
$d.Content.Find.Execute("11+43=54", $true, $false, $false, $false, $false, $true, 1, $false, "58-41=17", 2) | Out-Null
$d.Content.Find.Execute("90-14=76", $true, $false, $false, $false, $false, $true, 1, $false, "80-80=0", 2) | Out-Null
$d.Content.Find.Execute("52-43=9", $true, $false, $false, $false, $false, $true, 1, $false, "40+53=93", 2) | Out-Null
$d.Content.Find.Execute("62-30=32", $true, $false, $false, $false, $false, $true, 1, $false, "99-17=82", 2) | Out-Null
$d.Content.Find.Execute("82+5=87", $true, $false, $false, $false, $false, $true, 1, $false, "26-15=11", 2) | Out-Null
$d.Content.Find.Execute("51-14=37", $true, $false, $false, $false, $false, $true, 1, $false, "4+2=6", 2) | Out-Null
$d.Content.Find.Execute("44+41=85", $true, $false, $false, $false, $false, $true, 1, $false, "11+55=66", 2) | Out-Null
$d.Content.Find.Execute("16+69=85", $true, $false, $false, $false, $false, $true, 1, $false, "50-22=28", 2) | Out-Null
$d.Content.Find.Execute("5-2=3", $true, $false, $false, $false, $false, $true, 1, $false, "45-28=17", 2) | Out-Null
$d.Content.Find.Execute("54-26=28", $true, $false, $false, $false, $false, $true, 1, $false, "39+17=56", 2) | Out-Null
$d.Content.Find.Execute("8+13=21", $true, $false, $false, $false, $false, $true, 1, $false, "79-13=66", 2) | Out-Null
$d.Content.Find.Execute("54-24=30", $true, $false, $false, $false, $false, $true, 1, $false, "18+73=91", 2) | Out-Null
$d.Content.Find.Execute("87-32=55", $true, $false, $false, $false, $false, $true, 1, $false, "8-4=4", 2) | Out-Null
$d.Content.Find.Execute("33+10=43", $true, $false, $false, $false, $false, $true, 1, $false, "69-2=67", 2) | Out-Null
$d.Content.Find.Execute("96-63=33", $true, $false, $false, $false, $false, $true, 1, $false, "6-2=4", 2) | Out-Null
$d.Content.Find.Execute("33+18=51", $true, $false, $false, $false, $false, $true, 1, $false, "65+6=71", 2) | Out-Null
